$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = -20.912
$ws.Range("A13").Value = -22.005
$ws.Range("A16").Value = -20.86
$ws.Range("A18").Value = -21.868
$ws.Range("A20").Value = -21.883
